$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows of scenario results to the bottom of the sheet
# (rows 208-210), matching the existing data pattern (columns A=Scenario,
# B=Status, C=Browser).

$newRows = @(
    @("Create Nationality and Delete", "PASSED", "chrome"),
    @("Create Nationality and Delete", "PASSED", "chrome"),
    @("Login with valid username and password", "PASSED", "chrome")
)

$startRow = 208
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
}
